$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely; Excel shifts B->A, C->B, D->C for every row,
# which is exactly the layout change in the diff (the old index column A
# disappears and Sensibilizador/REGISTROS/Domicilios_sensibilizados move
# one column to the left).
$ws.Columns("A:A").Delete()
